# Added option for user to show map location for stores.
#
# This adds a new "lnglat" column (column N) containing the
# latitude,longitude coordinates for every store listed in the sheet, so
# that a downstream consumer (e.g. the Telegram bot) can show a map
# location for each outlet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 is the header ("lnglat"); rows 2-27 are the stores, in the same
# order as they already appear in column A of the sheet.
$lnglat = @{
    1  = "lnglat"
    2  = "1.3044719,103.7724654"
    3  = "1.307338,103.7726078"
    4  = "1.3045409,103.7727869"
    5  = "1.303794,103.7735167"
    6  = "1.306019,103.772678"
    7  = "1.3050106,103.7723947"
    8  = "1.3039258,103.7735858"
    9  = "1.3046285,103.7730182"
    10 = "1.3047189,103.7727242"
    11 = "1.3044706,103.7724575"
    12 = "1.3051092,103.7723276"
    13 = "1.3040592,103.7741032"
    14 = "1.305593,103.773083"
    15 = "1.3045187,103.7728417"
    16 = "1.3048207,103.7725693"
    17 = "1.305796,103.773008"
    18 = "1.3046285,103.7730182"
    19 = "1.3042717,103.7738946"
    20 = "1.3046387,103.7728153"
    21 = "1.3045756,103.7726986"
    22 = "1.3054322,103.7728657"
    23 = "1.3038022,103.7738266"
    24 = "1.3049764,103.7724652"
    25 = "1.3039101,103.7738303"
    26 = "1.3040203,103.7741394"
    27 = "1.3047292,103.7725536"
}

foreach ($row in $lnglat.Keys) {
    $ws.Cells.Item($row, 14).Value = $lnglat[$row]
}

# Reflect the newly added column as the cell the user is now looking at
# in the frozen (bottom-right) pane.
$ws.Range("N18").Select() | Out-Null
